# Add a new "2022-Q4" sheet (right after "总计") with fund-holding data,
# and update the "总计" (summary) sheet with a new leading row for 2022-Q4,
# cascading the previously-existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计" (i.e. before
#    the current second sheet, "2022-Q2"), matching the sheet order:
#    总计, 2022-Q4, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2
#    Duplicate "2022-Q2" (rather than adding a blank sheet) so the new
#    sheet inherits the same sheet-level setup (outline props, page
#    margins, etc.) as its siblings; its cell contents are overwritten
#    below.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2022-Q2")
$src.Copy($src, $null)
$q4 = $wb.Worksheets.Item("2022-Q2 (2)")
$q4.Name = "2022-Q4"

# Add the 3rd row (header + 2 data rows) by cloning the existing data
# row's formatting onto it.
$q4.Range("A2:H2").Copy()
$q4.Range("A3:H3").PasteSpecial(-4122)

# Header row text (identical wording/order to the other quarterly sheets)
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Data rows. The fund code (B) and the D/E/F/G figures are stored as
# plain *text* in the source workbook (e.g. "015921" keeps its leading
# zero). Assigning those verbatim would make Excel auto-coerce them to
# numbers, so momentarily mark the cell as Text (NumberFormat "@")
# while writing the value, then flip the style back to Normal/General
# so the saved cell ends up with no special formatting - exactly like
# the sibling quarterly sheets - while still keeping the text type.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$q4.Range("A2").Value = 0
Set-TextValue $q4.Range("B2") "015921"
$q4.Range("C2").Value = "申万菱信国证2000指数增强A"
Set-TextValue $q4.Range("D2") "0.21"
Set-TextValue $q4.Range("E2") "94.00"
Set-TextValue $q4.Range("F2") "0.51"
Set-TextValue $q4.Range("G2") "0.0011"
$q4.Range("H2").Value = 6

$q4.Range("A3").Value = 1
Set-TextValue $q4.Range("B3") "015922"
$q4.Range("C3").Value = "申万菱信国证2000指数增强C"
Set-TextValue $q4.Range("D3") "0.08"
Set-TextValue $q4.Range("E3") "94.00"
Set-TextValue $q4.Range("F3") "0.51"
Set-TextValue $q4.Range("G3") "0.0004"
$q4.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: push the existing quarterly rows
#    down by one and write a new row for "2022-Q4" at the top (row 2).
#    Column A is a plain 0-based row index, so it is simply rewritten
#    0..5 top to bottom, independent of which quarter moved where.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Make room for the new row by copying the existing index-column style
# onto the about-to-be-added row 7 before it holds any data.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

$rows = @(
    @{A=0; B="2022-Q4"; C=2; D=0},
    @{A=1; B="2022-Q2"; C=1; D=0.04},
    @{A=2; B="2022-Q1"; C=2; D=0},
    @{A=3; B="2021-Q4"; C=5; D=0.15},
    @{A=4; B="2021-Q3"; C=2; D=0},
    @{A=5; B="2021-Q2"; C=2; D=0}
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $total.Range("A$r").Value = $row.A
    $total.Range("B$r").Value = $row.B
    $total.Range("C$r").Value = $row.C
    $total.Range("D$r").Value = $row.D
}

Write-Output "2022-Q4 sheet added and 总计 table updated"
